$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.0070610687022900761
$ws.Range("C3").Value = 0.014122137404580152
$ws.Range("C4").Value = 0.0066793893129770982
$ws.Range("C5").Value = 0.010305343511450382

$ws.Range("D2").Value = "U"
$ws.Range("D3").Value = "U"
$ws.Range("D4").Value = "U"
$ws.Range("D5").Value = "U"

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1

$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = 2

$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2

$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 3
$ws.Range("G5").Value = 3

$ws.Range("D2:G5").Select()
